$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "도전! 한국 가요 마스터! (난이도 중)"
$ws.Range("B6").Value = "들려오는 곡의 '제목을' 입력해 주세요!"

$ws.Range("B12").Select()
